$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Datos actualizados a 9 de Abril de 2020 a las 21:22"

$ws.Cells.Item(4,2).Value = 456820
$ws.Cells.Item(4,3).Value = 21893
$ws.Cells.Item(4,5).Value = 416010
$ws.Cells.Item(4,7).Value = 1443
$ws.Cells.Item(4,8).Value = 16231

$ws.Cells.Item(19,2).Value = 13244
$ws.Cells.Item(19,3).Value = 302
$ws.Cells.Item(19,5).Value = 7709

$ws.Cells.Item(32,1).Value = "Peru"
$ws.Cells.Item(32,2).Value = 5256
$ws.Cells.Item(32,3).Value = 914
$ws.Cells.Item(32,4).Value = 1438
$ws.Cells.Item(32,5).Value = 3680
$ws.Cells.Item(32,6).Value = 124
$ws.Cells.Item(32,7).Value = 17
$ws.Cells.Item(32,8).Value = 138

$ws.Cells.Item(33,1).Value = "Rumania"
$ws.Cells.Item(33,2).Value = 5202
$ws.Cells.Item(33,3).Value = 441
$ws.Cells.Item(33,4).Value = 647
$ws.Cells.Item(33,5).Value = 4307
$ws.Cells.Item(33,6).Value = 178
$ws.Cells.Item(33,7).Value = 28
$ws.Cells.Item(33,8).Value = 248

$ws.Cells.Item(34,1).Value = "Ecuador"
$ws.Cells.Item(34,2).Value = 4965
$ws.Cells.Item(34,3).Value = 515
$ws.Cells.Item(34,4).Value = 339
$ws.Cells.Item(34,5).Value = 4354
$ws.Cells.Item(34,6).Value = 139
$ws.Cells.Item(34,7).Value = 30
$ws.Cells.Item(34,8).Value = 272

$ws.Cells.Item(35,1).Value = "Japon"
$ws.Cells.Item(35,2).Value = 4667
$ws.Cells.Item(35,3).Value = 0
$ws.Cells.Item(35,4).Value = 632
$ws.Cells.Item(35,5).Value = 3941
$ws.Cells.Item(35,6).Value = 99
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,8).Value = 94

$ws.Cells.Item(36,1).Value = "Pakistan"
$ws.Cells.Item(36,2).Value = 4489
$ws.Cells.Item(36,3).Value = 226
$ws.Cells.Item(36,4).Value = 572
$ws.Cells.Item(36,5).Value = 3852
$ws.Cells.Item(36,6).Value = 31
$ws.Cells.Item(36,7).Value = 4
$ws.Cells.Item(36,8).Value = 65

$ws.Cells.Item(56,1).Value = "Egipto"
$ws.Cells.Item(56,2).Value = 1699
$ws.Cells.Item(56,3).Value = 139
$ws.Cells.Item(56,4).Value = 348
$ws.Cells.Item(56,5).Value = 1233
$ws.Cells.Item(56,6).Value = 0
$ws.Cells.Item(56,7).Value = 15
$ws.Cells.Item(56,8).Value = 118

$ws.Cells.Item(57,1).Value = "Argelia"
$ws.Cells.Item(57,2).Value = 1666
$ws.Cells.Item(57,3).Value = 94
$ws.Cells.Item(57,4).Value = 347
$ws.Cells.Item(57,5).Value = 1084
$ws.Cells.Item(57,6).Value = 46
$ws.Cells.Item(57,7).Value = 30
$ws.Cells.Item(57,8).Value = 235

$ws.Cells.Item(58,1).Value = "Islandia"
$ws.Cells.Item(58,2).Value = 1648
$ws.Cells.Item(58,3).Value = 32
$ws.Cells.Item(58,4).Value = 688
$ws.Cells.Item(58,5).Value = 954
$ws.Cells.Item(58,6).Value = 11
$ws.Cells.Item(58,8).Value = 6

$ws.Cells.Item(87,1).Value = "Costa Rica"
$ws.Cells.Item(87,2).Value = 539
$ws.Cells.Item(87,3).Value = 37
$ws.Cells.Item(87,4).Value = 30
$ws.Cells.Item(87,5).Value = 506
$ws.Cells.Item(87,6).Value = 13
$ws.Cells.Item(87,7).Value = 0
$ws.Cells.Item(87,8).Value = 3

$ws.Cells.Item(88,1).Value = "Cuba"
$ws.Cells.Item(88,2).Value = 515
$ws.Cells.Item(88,3).Value = 58
$ws.Cells.Item(88,4).Value = 28
$ws.Cells.Item(88,5).Value = 472
$ws.Cells.Item(88,7).Value = 3
$ws.Cells.Item(88,8).Value = 15

$ws.Cells.Item(92,1).Value = "Costa de Marfil"
$ws.Cells.Item(92,2).Value = 444
$ws.Cells.Item(92,3).Value = 60
$ws.Cells.Item(92,4).Value = 52
$ws.Cells.Item(92,5).Value = 389
$ws.Cells.Item(92,7).Value = 0
$ws.Cells.Item(92,8).Value = 3

$ws.Cells.Item(93,1).Value = "Burkina Faso"
$ws.Cells.Item(93,2).Value = 443
$ws.Cells.Item(93,3).Value = 29
$ws.Cells.Item(93,4).Value = 146
$ws.Cells.Item(93,5).Value = 273
$ws.Cells.Item(93,6).Value = 0
$ws.Cells.Item(93,8).Value = 24

$ws.Cells.Item(94,1).Value = "Albania"
$ws.Cells.Item(94,2).Value = 409
$ws.Cells.Item(94,3).Value = 9
$ws.Cells.Item(94,4).Value = 165
$ws.Cells.Item(94,5).Value = 221
$ws.Cells.Item(94,6).Value = 7
$ws.Cells.Item(94,7).Value = 1
$ws.Cells.Item(94,8).Value = 23

$ws.Cells.Item(127,4).Value = 62
$ws.Cells.Item(127,5).Value = 57

$ws.Cells.Item(140,4).Value = 12
$ws.Cells.Item(140,5).Value = 47
